$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 5.794181456641729

# Row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 5.586269137925634
